$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on the Price/Volume columns so the updated
# numeric-looking strings (e.g. "1.000", "0.5090") are preserved exactly
# as text instead of being auto-converted to numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "25.819.54"
$ws.Range("D3").Value = "1.817.63"
$ws.Range("E3").Value = "  -2.97%  "
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "277.06"
$ws.Range("E5").Value = "  -7.98%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5090"
$ws.Range("E7").Value = "  -4.49%  "
$ws.Range("D8").Value = "0.3532"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("D9").Value = "44.62"
$ws.Range("E9").Value = "  -2.16%  "
$ws.Range("D10").Value = "0.06668"
$ws.Range("E10").Value = "  -7.09%  "
$ws.Range("D11").Value = "20.08"
$ws.Range("E11").Value = "  -7.18%  "
$ws.Range("D12").Value = "0.8298"
$ws.Range("E12").Value = "  -6.49%  "
$ws.Range("D13").Value = "0.07864"
$ws.Range("E13").Value = "  -3.72%  "
$ws.Range("D14").Value = "1.812.48"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "5.082"
$ws.Range("E15").Value = "  -3.90%  "
$ws.Range("D16").Value = "87.62"
$ws.Range("E16").Value = "  -6.16%  "
$ws.Range("D17").Value = "0.9996"
$ws.Range("E17").Value = "  -0.16%  "
$ws.Range("D18").Value = "14.15"
$ws.Range("E18").Value = "  -4.35%  "
$ws.Range("D19").Value = "0.000008055"
$ws.Range("E19").Value = "  -5.77%  "
$ws.Range("D20").Value = "1.000"
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("D21").Value = "25.866.87"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("D22").Value = "4.742"
$ws.Range("E22").Value = "  -4.85%  "
$ws.Range("D23").Value = "10.02"
$ws.Range("E23").Value = "  -6.20%  "
$ws.Range("D24").Value = "6.095"
$ws.Range("E24").Value = "  -4.64%  "
$ws.Range("D25").Value = "2.204"
$ws.Range("E25").Value = "  -3.14%  "
$ws.Range("D26").Value = "141.26"
$ws.Range("E26").Value = "  -3.52%  "
$ws.Range("D27").Value = "1.675"
$ws.Range("E27").Value = "  -3.63%  "
$ws.Range("D28").Value = "17.10"
$ws.Range("E28").Value = "  -5.18%  "
$ws.Range("D29").Value = "109.50"
$ws.Range("E29").Value = "  -3.81%  "
$ws.Range("D30").Value = "4.359"
$ws.Range("E30").Value = "  -7.84%  "
$ws.Range("D31").Value = "4.247"
$ws.Range("E31").Value = "  -8.07%  "
$ws.Range("D32").Value = "0.08801"
$ws.Range("D33").Value = "0.04897"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("D34").Value = "0.7314"
$ws.Range("E34").Value = "  -9.84%  "
$ws.Range("D35").Value = "1.142"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("E36").Value = "  -3.07%  "
$ws.Range("D37").Value = "0.9997"
$ws.Range("D38").Value = "3.135"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").Value = "2.391"
$ws.Range("E39").Value = "  -7.78%  "
$ws.Range("D40").Value = "0.5206"
$ws.Range("E40").Value = "  -14.02%  "
$ws.Range("D41").Value = "0.01854"
$ws.Range("E41").Value = "  -5.02%  "
$ws.Range("D42").Value = "0.9632"
$ws.Range("E42").Value = "  -10.05%  "
$ws.Range("D43").Value = "6.216"
$ws.Range("E43").Value = "  -5.48%  "
$ws.Range("D44").Value = "110.82"
$ws.Range("E44").Value = "  -3.36%  "
$ws.Range("D45").Value = "8.042"
$ws.Range("E45").Value = "  -9.42%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("D47").Value = "0.4590"
$ws.Range("E47").Value = "  -10.80%  "
$ws.Range("D48").Value = "0.1367"
$ws.Range("E48").Value = "  -8.49%  "
$ws.Range("D49").Value = "36.68"
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("D50").Value = "9.296"
$ws.Range("E50").Value = "  -6.10%  "
$ws.Range("D51").Value = "1.505"
$ws.Range("E51").Value = "  -7.88%  "

# Restore the default (unstyled) cell style now that the values are
# committed as text, matching the original workbook formatting.
$ws.Range("D2:E51").Style = "Normal"
